# Updates the cryptos list values/percentages on Sheet1, matching the
# upstream GitHub Actions data refresh described in the commit message.
#
# Note: several "Price" column values look like plain numbers to Excel
# (e.g. "253.32") and would otherwise be auto-converted from text to a
# numeric value, dropping significant trailing zeros. We prefix those
# with a literal leading apostrophe (the standard Excel "force text"
# convention) so they are stored as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.135.90'
$ws.Range("D3").Value = '1.901.51'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").Value = '''253.32'
$ws.Range("E5").Value = '  +3.12%  '
$ws.Range("D6").Value = '''0.692'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").Value = '''41.42'
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("E9").Value = '  +2.95%  '
$ws.Range("D10").Value = '''52.77'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  +3.85%  '
$ws.Range("D12").Value = '''0.0983'
$ws.Range("E12").Value = '  -1.21%  '
# Row 13/14: coin identities swap (Chainlink <-> WrappedliquidstakedEther2.0)
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.178.05'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''13.01'
$ws.Range("E14").Value = '  +4.73%  '
$ws.Range("D15").Value = '''0.734'
$ws.Range("E15").Value = '  +3.88%  '
$ws.Range("D16").Value = '''4.98'
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").Value = '1.900.00'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '35.136.23'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '''73.56'
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '''243.35'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").Value = '''12.97'
$ws.Range("E22").Value = '  +2.74%  '
$ws.Range("E23").Value = '  +5.10%  '
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").Value = '''2.43'
$ws.Range("E25").Value = '  +5.04%  '
$ws.Range("D26").Value = '''2.28'
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("D27").Value = '''166.94'
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = '''8.56'
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").Value = '''18.52'
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("D30").Value = '''0.130'
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").Value = '4.128.84'
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").Value = '''2.05'
$ws.Range("E32").Value = '  +12.80%  '
$ws.Range("E33").Value = '  +6.90%  '
$ws.Range("E34").Value = '  +4.00%  '
$ws.Range("E35").Value = '  +6.72%  '
$ws.Range("E36").Value = '  +2.64%  '
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("D38").Value = '''0.857'
$ws.Range("E38").Value = '  -7.63%  '
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("D40").Value = '''102.93'
$ws.Range("E40").Value = '  +14.95%  '
$ws.Range("D41").Value = '''17.40'
$ws.Range("E41").Value = '  +8.84%  '
$ws.Range("E42").Value = '  +2.55%  '
$ws.Range("E43").Value = '  +0.90%  '
$ws.Range("E44").Value = '  -2.07%  '
$ws.Range("D45").Value = '''2.42'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = '1.320.28'
$ws.Range("E46").Value = '  -2.29%  '
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("D50").Value = '''12.10'
$ws.Range("E50").Value = '  -4.41%  '
# Row 51: coin identity changes from Cronos to MultiversX
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''43.48'
$ws.Range("E51").Value = '  -5.62%  '
